# Update the "Cost / project" (column G) values for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G12").Value = 12320000
$ws.Range("G29").Value = 260279.126780675
$ws.Range("G30").Value = 395666.69261350599
$ws.Range("G32").Value = 464784.15496549098
$ws.Range("G33").Value = 6051489.6976506999

# Match the author's final selection (G29:G33, active cell G29).
[void]$ws.Range("G29:G33").Select()
